# ALZ Policy Assignments v2.xlsx -- "Policy Refresh Q1 FY24" edit
#
# Applies the content changes described by the commit diff:
#   - Row 15 (Key Vault guardrails initiative) fully repopulated + new hyperlink
#   - Row 26 / Row 33 column G: "Audit" -> "Deny"
#   - Row 43 / Row 44 column A: shortened scope labels
#   - Selection / scroll position moved to reflect where the author was working
#
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "ALZ Policy Assignments 03CY23"

# ---------------------------------------------------------------------------
# Row 15 - was a placeholder "N/A" row, now the Key Vault guardrails entry
# ---------------------------------------------------------------------------
$ws.Rows.Item(15).RowHeight = 43.2

$ws.Cells.Item(15, 1).Value = "Platform"
$ws.Cells.Item(15, 2).Value = "Enforce recommendded guardrails for Azure Key Vault"
$ws.Cells.Item(15, 3).Value = "Enforce recommendded guardrails for Azure Key Vault"
$ws.Cells.Item(15, 4).Value = "Initiative"
$ws.Cells.Item(15, 5).Value = "Custom"
$ws.Cells.Item(15, 6).Value = "This initiative assignment enables recommended ALZ guardrails for Azure Key Vault."
$ws.Cells.Item(15, 7).Value = "Deny, Audit"
$ws.Cells.Item(15, 8).Value = "ENFORCE-GuardrailsKeyVaultPolicyAssignment.json"
$ws.Cells.Item(15, 9).Value = "Enforce recommended guardrails for Azure Key Vault (azadvertizer.net)"
$ws.Cells.Item(15, 10).Value = 45124

# Turn I15 into a hyperlink pointing at the azadvertizer page, same as the
# other policy rows.
$ws.Hyperlinks.Add($ws.Cells.Item(15, 9), "https://www.azadvertizer.net/azpolicyinitiativesadvertizer/Enforce-Guardrails-KeyVault.html", "", "", "Enforce recommended guardrails for Azure Key Vault (azadvertizer.net)")

# Hyperlinks.Add stamps a brand-new cell style; restore the normal
# "hyperlink, wrap text, no vertical-top" look used by the other links
# (e.g. I30) by copying its format over.
$ws.Range("I30").Copy()
$ws.Range("I15").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Rows 26 & 33 - Effect column: Audit -> Deny
# ---------------------------------------------------------------------------
$ws.Cells.Item(26, 7).Value = "Deny"
$ws.Cells.Item(33, 7).Value = "Deny"

# ---------------------------------------------------------------------------
# Rows 43 & 44 - shorten the Assignment Scope (MG) labels
# ---------------------------------------------------------------------------
$ws.Cells.Item(43, 1).Value = "Decommissioned"
$ws.Cells.Item(44, 1).Value = "Sandbox"

# ---------------------------------------------------------------------------
# Sheet view - move the scroll position / active selection
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I28").Select()
